# them chuc nang chup anh man hinh
# Append new "ban tin da gui" (sent message log) rows describing newly
# received alarm notifications, and touch a handful of pre-existing
# timestamp cells whose sub-millisecond precision was refreshed by the
# source system on re-export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the sub-millisecond jitter on a few already-logged rows ---
$ws.Cells.Item(116, 1).Value = 45754.73332535879
$ws.Cells.Item(117, 1).Value = 45754.73333599537
$ws.Cells.Item(118, 1).Value = 45754.73334625
$ws.Cells.Item(119, 1).Value = 45754.73335608796
$ws.Cells.Item(120, 1).Value = 45754.73336722222

# --- Append the newly received notifications ---
function Add-LogRow {
    param($Row, $SentAt, $NeName, $AlarmContent, $Status)

    $cellA = $ws.Cells.Item($Row, 1)
    $cellA.Value = $SentAt
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($Row, 2).Value = $NeName
    $ws.Cells.Item($Row, 3).Value = $AlarmContent
    $ws.Cells.Item($Row, 4).Value = $Status
}

Add-LogRow 121 45755.34974707176 "UL_STY078M_HNI" "POWER_AC_EAS" "Thành công"
Add-LogRow 122 45755.36286347222 "UL_TTT038M_HNI" "POWER_AC_EAS" "Thành công"
Add-LogRow 123 45755.41738707176 "SR_TTT014M_HNI" "POWER_AC_EAS" "Thành công"
Add-LogRow 124 45755.42811016204 "2G_BVI015M_HNI" "SITE_OOS" "Thành công"
Add-LogRow 125 45755.42812280093 "3G_BVI015M_HNI" "SITE_OOS" "Thành công"
Add-LogRow 126 45755.43242371528 "UL_TTT111M_HNI" "POWER_AC_EAS" "Thành công"
Add-LogRow 127 45755.44526734954 "3G_STY030M_HNI" "SITE_OOS" "Thành công"
Add-LogRow 128 45755.44527988426 "2G_STY030M_HNI" "SITE_OOS" "Thành công"
Add-LogRow 129 45757.91495641533 "4G_STY022M_HNI" "CELL_OOS, HW_OUTDOOR, HW_TTS_MUC_TRAM" "Thành công"
